$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the "Model holdings provided as of ..." disclaimer text in A9.
$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-14 for illustrative purposes only and are subject to change."

# Update the Weight (D) and Percent Change (E) values for rows 2-6.
$ws.Range("D2").Value = 0.2572267077861345
$ws.Range("E2").Value = 0.01197191205249215

$ws.Range("D3").Value = 0.2548315560650698
$ws.Range("E3").Value = 0.01582193617591843

$ws.Range("D4").Value = 0.2452246495492019
$ws.Range("E4").Value = 0.01384317521781231

$ws.Range("D5").Value = 0.2427170865995939
$ws.Range("E5").Value = 0.01970899470899479

$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 0.01528982170584947

# Restore the sheet protection that was in place before the edit.
$ws.Protect()
